$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-26 08:43:24"

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-26 08:43:20"

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-26 08:43:24"

# --- Column width adjustments (status column widened to fit "Ready for handoff") ---
# Target stored width (per XML diff) is 17.2159881591797 characters; the closest
# value this engine's ColumnWidth rounding can reach is 17.1666... , which is
# produced by any input in roughly [16.25, 16.41].
$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25
$wsZhCn.Columns.Item(3).ColumnWidth = 16.25
$wsDeDe.Columns.Item(3).ColumnWidth = 16.25
